# Regenerate save_data column G ("K" = strikeouts) with freshly scraped
# values, replacing the stale "Strike#" derived numbers.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (1-based, matching the worksheet) -> new value for column G ("K")
$gUpdates = @(
    @{ Row = 2; Value = 0 },
    @{ Row = 3; Value = 2 },
    @{ Row = 4; Value = 0 },
    @{ Row = 5; Value = 2 },
    @{ Row = 6; Value = 2 },
    @{ Row = 8; Value = 0 },
    @{ Row = 9; Value = 2 },
    @{ Row = 10; Value = 2 },
    @{ Row = 11; Value = 1 },
    @{ Row = 12; Value = 1 },
    @{ Row = 13; Value = 2 },
    @{ Row = 14; Value = 0 },
    @{ Row = 15; Value = 2 },
    @{ Row = 16; Value = 1 },
    @{ Row = 17; Value = 1 },
    @{ Row = 18; Value = 0 },
    @{ Row = 19; Value = 1 },
    @{ Row = 20; Value = 2 },
    @{ Row = 21; Value = 0 },
    @{ Row = 22; Value = 2 },
    @{ Row = 23; Value = 1 },
    @{ Row = 24; Value = 1 },
    @{ Row = 25; Value = 1 },
    @{ Row = 27; Value = 0 },
    @{ Row = 28; Value = 2 },
    @{ Row = 29; Value = 2 },
    @{ Row = 30; Value = 1 },
    @{ Row = 31; Value = 1 },
    @{ Row = 32; Value = 0 },
    @{ Row = 33; Value = 1 },
    @{ Row = 34; Value = 1 },
    @{ Row = 35; Value = 2 },
    @{ Row = 36; Value = 1 },
    @{ Row = 37; Value = 1 },
    @{ Row = 38; Value = 0 },
    @{ Row = 39; Value = 2 },
    @{ Row = 40; Value = 0 },
    @{ Row = 41; Value = 2 },
    @{ Row = 42; Value = 2 },
    @{ Row = 43; Value = 0 },
    @{ Row = 44; Value = 1 },
    @{ Row = 45; Value = 1 },
    @{ Row = 46; Value = 2 },
    @{ Row = 47; Value = 2 },
    @{ Row = 48; Value = 1 },
    @{ Row = 49; Value = 0 },
    @{ Row = 50; Value = 1 },
    @{ Row = 51; Value = 3 },
    @{ Row = 52; Value = 1 },
    @{ Row = 53; Value = 0 },
    @{ Row = 54; Value = 2 },
    @{ Row = 55; Value = 1 },
    @{ Row = 56; Value = 1 },
    @{ Row = 57; Value = 2 },
    @{ Row = 58; Value = 1 },
    @{ Row = 59; Value = 1 },
    @{ Row = 60; Value = 0 },
    @{ Row = 61; Value = 2 },
    @{ Row = 62; Value = 0 },
    @{ Row = 63; Value = 1 },
    @{ Row = 64; Value = 1 },
    @{ Row = 65; Value = 2 },
    @{ Row = 66; Value = 1 },
    @{ Row = 67; Value = 2 },
    @{ Row = 68; Value = 2 },
    @{ Row = 69; Value = 2 },
    @{ Row = 70; Value = 0 },
    @{ Row = 71; Value = 1 },
    @{ Row = 72; Value = 0 },
    @{ Row = 73; Value = 2 },
    @{ Row = 74; Value = 1 },
    @{ Row = 75; Value = 1 },
    @{ Row = 76; Value = 1 },
    @{ Row = 77; Value = 1 },
    @{ Row = 78; Value = 1 },
    @{ Row = 79; Value = 1 }
)

foreach ($update in $gUpdates) {
    $ws.Cells.Item($update.Row, 7).Value = $update.Value
}

$wb.Save()
